$wb = $excel.ActiveWorkbook

# --- Sheet "Folha Ponto": add two new time entries (row 25 & 26) ---
$ws1 = $wb.Worksheets.Item("Folha Ponto")

# Row 25: 10:30 - 13:15, project SITS, activity "Trabalhando em novo level design"
$ws1.Range("B25").Value = 44906
$ws1.Range("C25").Value = 0.4375
$ws1.Range("D25").Value = 0.55208333333333337
$ws1.Range("F25").Value = "SITS"
$ws1.Range("G25").Value = "Trabalhando em novo level design"

# Row 26: 14:00 - 17:00, project SITS, activity "Trabalhando em novo level design"
$ws1.Range("B26").Value = 44906
$ws1.Range("C26").Value = 0.58333333333333337
$ws1.Range("D26").Value = 0.70833333333333337
$ws1.Range("F26").Value = "SITS"

$g26 = $ws1.Range("G26")
$g26.ClearFormats() | Out-Null
$g26.Value = "Trabalhando em novo level design"

# --- Sheet "Totais": selection only moves; totals recalc automatically ---
$ws2 = $wb.Worksheets.Item("Totais")
$ws2.Activate() | Out-Null
$ws2.Range("H6").Select() | Out-Null

# Leave "Folha Ponto" as the active sheet/selection when done
$ws1.Activate() | Out-Null
$ws1.Range("D27").Select() | Out-Null
